$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'baseball pants high knee'
$ws.Range("A2").Value = 'knee pack'
$ws.Range("A3").Value = 'girls compression pants'
$ws.Range("A4").Value = 'basketball gear for boys'
$ws.Range("A5").Value = 'gym knee pads'
$ws.Range("A6").Value = 'extra large knee pads'
$ws.Range("A7").Value = 'tall leggings'
$ws.Range("A8").Value = 'large knee pads'
$ws.Range("A9").Value = 'football shock pads'
$ws.Range("A10").Value = 'fitness compression pants men'
$ws.Range("A11").Value = 'spandex leggings'
$ws.Range("A12").Value = 'knee pad hiking'
$ws.Range("A13").Value = 'workout legging for men'
$ws.Range("A14").Value = 'compression legs men'
$ws.Range("A15").Value = 'pants for men'
$ws.Range("A16").Value = 'mens baseball clothing'
$ws.Range("A17").Value = 'capri legging pants'
$ws.Range("A18").Value = 'compression pants knee length'
$ws.Range("A19").Value = 'knee pad for pain'
$ws.Range("A20").Value = 'camo knee pads for basketball'
$ws.Range("A21").Value = 'mcdavid basketball tights knee pads'
$ws.Range("A22").Value = 'skateboard knee pads youth'
$ws.Range("A23").Value = 'kids basketball knee pads'
$ws.Range("A24").Value = 'youth basketball pants tearaway'
$ws.Range("A25").Value = 'kneee pads for basketball'
$ws.Range("A26").Value = 'commpression pants for basketball'
$ws.Range("A27").Value = 'supreme basketball leggings'
$ws.Range("A28").Value = 'volleyball knee pads youth girls'
$ws.Range("A29").Value = 'adidas knee pads basketball'
$ws.Range("A30").Value = 'basketball tights men'
$ws.Range("A31").Value = 'compression pants men adidas'
$ws.Range("A32").Value = 'pant with knee pads'
$ws.Range("A33").Value = 'knee pads tights'
$ws.Range("A34").Value = 'compression pants with knee'
$ws.Range("A35").Value = 'women basketball pants'
$ws.Range("A36").Value = 'basketball knee pads for youth kids'
$ws.Range("A37").Value = 'pantalones con rodilleras'
$ws.Range("A38").Value = 'caterpillar knee pad pants'
$ws.Range("A39").Value = 'thick mens leggings'
$ws.Range("A40").Value = 'capri leggings men'
$ws.Range("A41").Value = 'compression leggings men basketball'
$ws.Range("A42").Value = 'capri tights men'
$ws.Range("A43").Value = 'mens compression knee pads'
$ws.Range("A44").Value = 'legging men'
$ws.Range("A45").Value = 'basketball knee pads youth'
$ws.Range("A46").Value = 'hex pads basketball knee'
$ws.Range("A47").Value = 'baseball knee high pants'
$ws.Range("A48").Value = 'little boys baseball pants'
$ws.Range("A49").Value = 'boys baseball pants'
$ws.Range("A50").Value = 'knee pad compression'
$ws.Range("A51").Value = 'black compression tights'
$ws.Range("A52").Value = 'capris men'
$ws.Range("A53").Value = 'hip pads for volleyball'
$ws.Range("A54").Value = 'soccer pants men'
$ws.Range("A55").Value = 'boys knee pads'
$ws.Range("A56").Value = 'capris pants men'
$ws.Range("A57").Value = 'knee pad baseball'
$ws.Range("A58").Value = 'volleyball kneepads youth'
$ws.Range("A59").Value = 'compression leggings for boys'
$ws.Range("A60").Value = 'knee pads lightweight'
$ws.Range("A61").Value = 'knee pads sliding'
$ws.Range("A62").Value = 'capri leggings mesh'
$ws.Range("A63").Value = 'men running pants'
$ws.Range("A64").Value = 'knee pad black'
$ws.Range("A65").Value = 'baseball pants youth boys'
$ws.Range("A66").Value = 'patella protector'
$ws.Range("A67").Value = 'leggings pack'
$ws.Range("A68").Value = 'knee pads for workout'
$ws.Range("A69").Value = 'capri pants boys'
$ws.Range("A70").Value = 'gym pad men'
$ws.Range("A71").Value = 'mens athletic pants'
$ws.Range("A72").Value = 'leggings youth'
$ws.Range("A73").Value = 'knee pads protection'
$ws.Range("A74").Value = 'knee pad softball'
$ws.Range("A75").Value = 'knee pads fitness'
$ws.Range("A76").Value = 'leggings youth boys'
$ws.Range("A77").Value = 'gel knee pad'
$ws.Range("A78").Value = 'basketballs for boys'
$ws.Range("A79").Value = 'women athletic leggings'
$ws.Range("A80").Value = 'knee pads insert'
$ws.Range("A81").Value = 'capri leggings'
$ws.Range("A82").Value = 'elbow pads basketball youth'
$ws.Range("A83").Value = 'nike compression tights'
$ws.Range("A84").Value = 'nike compression leggings men'
$ws.Range("A85").Value = 'gloves and knee pads'
$ws.Range("A86").Value = 'black leggings running'
$ws.Range("A87").Value = 'lupo compression leggings'
$ws.Range("A88").Value = 'compression leggings 30'
$ws.Range("A89").Value = 'compression leggings black'
$ws.Range("A90").Value = 'compression leggings circulation'
$ws.Range("A91").Value = 'compression leggings girls'
$ws.Range("A92").Value = 'compression leggings men under armour'
$ws.Range("A93").Value = 'compression leggings men white'
$ws.Range("A94").Value = 'compression leggings men nike'
$ws.Range("A95").Value = 'compression leggings running'
$ws.Range("A96").Value = 'compression leggings xxl'
$ws.Range("A97").Value = 'running capri pants'
$ws.Range("A98").Value = 'knee basketball pads'
$ws.Range("A99").Value = 'basketball knee pads youth'
$ws.Range("A100").Value = 'basketball knee pads kids youth'
